# Auto-generated Excel COM-interop script
# Applies the cryptos list update (values + the RenderToken/Aave row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '58.816.26'
$ws.Range('E2').Value = '  +1.44%  '
$ws.Range('D3').Value = '2.501.11'
$ws.Range('E3').Value = '  +1.42%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '533.52'
$ws.Range('E5').Value = '  +3.91%  '
$ws.Range('D6').Value = '133.76'
$ws.Range('E6').Value = '  +2.66%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.35%  '
$ws.Range('D8').Value = '0.571'
$ws.Range('E8').Value = '  +3.61%  '
$ws.Range('D9').Value = '2.504.38'
$ws.Range('E9').Value = '  +0.63%  '
$ws.Range('D10').Value = '0.0992'
$ws.Range('E10').Value = '  +2.81%  '
$ws.Range('E11').Value = '  -2.54%  '
$ws.Range('D12').Value = '5.20'
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('E13').Value = '  +0.07%  '
$ws.Range('D14').Value = '2.936.10'
$ws.Range('E14').Value = '  +1.54%  '
$ws.Range('D15').Value = '58.647.30'
$ws.Range('E15').Value = '  +1.31%  '
$ws.Range('D16').Value = '22.28'
$ws.Range('E16').Value = '  +1.47%  '
$ws.Range('E17').Value = '  +1.66%  '
$ws.Range('D18').Value = '2.498.40'
$ws.Range('E18').Value = '  +1.14%  '
$ws.Range('D19').Value = '10.59'
$ws.Range('E19').Value = '  +0.20%  '
$ws.Range('E20').Value = '  +2.80%  '
$ws.Range('D21').Value = '320.61'
$ws.Range('E21').Value = '  +0.62%  '
$ws.Range('D22').Value = '6.21'
$ws.Range('E22').Value = '  +4.61%  '
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').Value = '66.00'
$ws.Range('E24').Value = '  +4.79%  '
$ws.Range('D25').Value = '0.407'
$ws.Range('E25').Value = '  +1.60%  '
$ws.Range('D26').Value = '0.992'
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('E27').Value = '  -0.63%  '
$ws.Range('E28').Value = '  +2.47%  '
$ws.Range('D29').Value = '173.14'
$ws.Range('E29').Value = '  +2.18%  '
$ws.Range('D30').Value = '0.0₃0756'
$ws.Range('E30').Value = '  +3.05%  '
$ws.Range('E31').Value = '  +3.66%  '
$ws.Range('E32').Value = '  +1.17%  '
$ws.Range('E33').Value = '  +0.55%  '
$ws.Range('E34').Value = '  +0.16%  '
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').Value = '  +0.38%  '
$ws.Range('E36').Value = '  +1.42%  '
$ws.Range('E37').Value = '  -3.13%  '
$ws.Range('D38').Value = '3.95'
$ws.Range('E38').Value = '  +1.03%  '
$ws.Range('D39').Value = '0.836'
$ws.Range('E39').Value = '  +9.22%  '
$ws.Range('D40').Value = '1.51'
$ws.Range('E40').Value = '  +3.56%  '
$ws.Range('D41').Value = '36.35'
$ws.Range('E41').Value = '  -0.76%  '
$ws.Range('D42').Value = '3.47'
$ws.Range('E42').Value = '  +2.56%  '
$ws.Range('D43').Value = '275.06'
$ws.Range('E43').Value = '  +1.00%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = '131.16'
$ws.Range('E44').Value = '  +8.24%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = '5.00'
$ws.Range('E45').Value = '  -0.69%  '
$ws.Range('D46').Value = '0.592'
$ws.Range('E46').Value = '  +0.70%  '
$ws.Range('D47').Value = '0.0933'
$ws.Range('E47').Value = '  +2.04%  '
$ws.Range('D48').Value = '0.0509'
$ws.Range('E48').Value = '  +4.18%  '
$ws.Range('E49').Value = '  +3.00%  '
$ws.Range('D50').Value = '16.78'
$ws.Range('E50').Value = '  +0.56%  '
$ws.Range('D51').Value = '1.754.36'
$ws.Range('E51').Value = '  +2.20%  '
